# Generate Report for Handback
# Update the handoff/handback datetime stamps recorded for the zh-cn and
# de-de localization rows to reflect the latest report generation run.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 2 holds the 81566471-...zh-cn.xlf entry
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-11 10:55:02"
$wsZhCn.Range("H2").Value = "2016-03-11 10:55:18"

# de-de sheet: row 2 holds the 81566471-...de-de.xlf entry
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-11 10:55:06"
$wsDeDe.Range("H2").Value = "2016-03-11 10:55:24"
